$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Résultat Actuel" column header (C1) to "Réalisation"
$ws.Range("C1").Value = "Réalisation"

# Move the selection to C30 (matches the saved selection state in the file)
$ws.Range("C30").Select()
